$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the empty paragraph that
#    currently sits between the "Echter, ... thee wil." paragraph and the
#    "Als je eenmaal ..." paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$wordmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Helper: append a brand-new paragraph (built from raw WordML) right after
# the current final paragraph of the document, without disturbing the
# paragraph that is already last.
function Add-ParagraphXml($xml) {
    $d2 = $word.ActiveDocument
    $lastPara = $d2.Paragraphs.Item($d2.Paragraphs.Count)
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d2.Paragraphs.Item($d2.Paragraphs.Count)
    [void]$newPara.Range.InsertXML($xml)
}

# 2. New blank paragraph (<w:p/>) right after the existing trailing blank
#    paragraph.
Add-ParagraphXml("<w:p $wordmlNs></w:p>")

# 3. "TODO:" paragraph, now carrying the relocated "_GoBack" bookmark.
Add-ParagraphXml(
  "<w:p $wordmlNs>" +
    "<w:r><w:t>TODO:</w:t></w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
    "<w:bookmarkEnd w:id=`"0`"/>" +
  "</w:p>"
)

# 4. "Wat ga ik opleveren en hoe ga ik het doen?" paragraph.
Add-ParagraphXml(
  "<w:p $wordmlNs>" +
    "<w:r><w:t>Wat ga ik opleveren en hoe ga ik het doen?</w:t></w:r>" +
  "</w:p>"
)

# 5. "Feedpulse opleveren volgende week" paragraph.
Add-ParagraphXml(
  "<w:p $wordmlNs>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Feedpulse</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> opleveren volgende week</w:t></w:r>" +
  "</w:p>"
)

# 6. "Bartosz contact opnemen CI, frontend testing enz.: hoe aantonen?"
#    paragraph.
Add-ParagraphXml(
  "<w:p $wordmlNs>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Bartosz</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> contact opnemen CI, </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>frontend</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>testing</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> en</w:t></w:r>" +
    "<w:r><w:t>z.: hoe aantonen?</w:t></w:r>" +
  "</w:p>"
)

Write-Output "All edits applied"
